# Finished revamping RQ1, added RQ1 values for new papers, finished drawing RQ1 figures
#
# - Re-assert the D2:D46 shared formula (shrinking it from D2:D49) so the
#   "C/B" ratio column no longer covers the now-removed D48 formula cell.
# - Fully remove the now-stray D48 formula cell (row 48 keeps its other
#   values; only the ratio formula goes away).
# - Append five new task rows (51-55) with their Page/Items counts, one of
#   which is computed with a small arithmetic formula.
# - Update the active sheet selection to reflect where the new rows were
#   added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-share the ratio formula over its (now shorter) range, then drop the
# leftover formula cell at D48 completely (no value, no formula there).
$ws.Range("D2:D46").Formula = "=C2/B2"
$ws.Range("D48").Clear()

# New rows for RQ1 follow-up work.
$ws.Range("A51").Value = "Work on the value of RQ1"
$ws.Range("B51").Value = 1
$ws.Range("C51").Formula = "=60*3+17"

$ws.Range("A52").Value = 'Update the technical problems, remove "remove trusted third party" and other generic values'
$ws.Range("B52").Value = 1
$ws.Range("C52").Value = 150

$ws.Range("A53").Value = "Update the improvement objectives to add proper new functionality"
$ws.Range("B53").Value = 1
$ws.Range("C53").Value = 12

$ws.Range("A54").Value = "Finish adding values for RQ1"
$ws.Range("B54").Value = 1
$ws.Range("C54").Value = 18

$ws.Range("A55").Value = "Finish drawing new figures for RQ1"
$ws.Range("B55").Value = 1
$ws.Range("C55").Value = 60

# Match the author's final selection/scroll position on the sheet.
$ws.Range("C51:C55").Select()
